# Natmi following Dr Hou advice
# Update Ligand-expressing cells (E) and Receptor-expressing cells (K) counts
# from 1 to 3, and recompute the dependent expression/specificity metrics
# (G, H, I, J, M, N, O, P, Q, R, S, T) for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989; K=3; M=10.000565; N=30.001695; O=0.6316353758144477; P=0.6316353758144477; Q=91.57476701535668; R=824.17290313821; S=0.6123512090853941; T=0.6123512090853941 }
    3 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989; K=3; M=4.264793333333333; N=12.79438; O=0.2693642149089528; P=0.2693642149089528; Q=39.05253911840445; R=351.47285206564; S=0.2611403809850739; T=0.2611403809850738 }
    4 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989; K=3; M=1.567455; N=4.702364999999999; O=0.09900040927659938; P=0.09900040927659938; Q=14.35312169183; R=129.17809522647; S=0.09597787369383093; T=0.09597787369383091 }
    5 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109; K=3; M=10.000565; N=30.001695; O=0.6316353758144477; P=0.6316353758144477; Q=2.883872929615; R=25.954856366535; S=0.01928416672905367; T=0.01928416672905367 }
    6 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109; K=3; M=4.264793333333333; N=12.79438; O=0.2693642149089528; P=0.2693642149089528; Q=1.229842718326667; R=11.06858446494; S=0.008223833923878958; T=0.008223833923878958 }
    7 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109; K=3; M=1.567455; N=4.702364999999999; O=0.09900040927659938; P=0.09900040927659938; Q=0.4520085658049999; R=4.068077092245; S=0.003022535582768455; T=0.003022535582768455 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
